$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = "'27.267.07"
$ws.Range('E2').Value = "'  +0.43%  "
$ws.Range('D3').Value = "'1.773.15"
$ws.Range('E3').Value = "'  +3.64%  "
$ws.Range('E4').Value = "'  +0.09%  "
$ws.Range('D5').Value = "'313.55"
$ws.Range('E5').Value = "'  +1.91%  "
$ws.Range('E6').Value = "'  +0.10%  "
$ws.Range('D7').Value = "'0.5267"
$ws.Range('E7').Value = "'  +9.95%  "
$ws.Range('E8').Value = "'  +6.12%  "
$ws.Range('D9').Value = "'42.73"
$ws.Range('E9').Value = "'  +1.86%  "
$ws.Range('E10').Value = "'  +1.08%  "
$ws.Range('E11').Value = "'  +3.91%  "
$ws.Range('E12').Value = "'  +0.15%  "
$ws.Range('D13').Value = "'20.45"
$ws.Range('E13').Value = "'  +2.87%  "
$ws.Range('D14').Value = "'6.058"
$ws.Range('E14').Value = "'  +3.49%  "
$ws.Range('D15').Value = "'1.767.09"
$ws.Range('E15').Value = "'  +3.45%  "
$ws.Range('D16').Value = "'6.946"
$ws.Range('E16').Value = "'  +1.51%  "
$ws.Range('D17').Value = "'88.69"
$ws.Range('E17').Value = "'  -0.53%  "
$ws.Range('E18').Value = "'  +0.64%  "
$ws.Range('E19').Value = "'  +1.17%  "
$ws.Range('D20').Value = "'1.001"
$ws.Range('D21').Value = "'16.73"
$ws.Range('E21').Value = "'  +1.58%  "
$ws.Range('E22').Value = "'  +4.14%  "
$ws.Range('D23').Value = "'27.350.17"
$ws.Range('E23').Value = "'  +0.61%  "
$ws.Range('D24').Value = "'11.25"
$ws.Range('E24').Value = "'  +3.40%  "
$ws.Range('E25').Value = "'  +0.00%  "
$ws.Range('D26').Value = "'154.55"
$ws.Range('E26').Value = "'  +0.20%  "
$ws.Range('E27').Value = "'  +2.18%  "
$ws.Range('D28').Value = "'2.328"
$ws.Range('E28').Value = "'  +11.90%  "
$ws.Range('D29').Value = "'1.968.98"
$ws.Range('E29').Value = "'  +3.48%  "
$ws.Range('D30').Value = "'120.93"
$ws.Range('E30').Value = "'  +1.15%  "
$ws.Range('D31').Value = "'1.060"
$ws.Range('E31').Value = "'  +5.05%  "
$ws.Range('D32').Value = "'0.09755"
$ws.Range('E32').Value = "'  +5.12%  "
$ws.Range('D33').Value = "'5.563"
$ws.Range('E33').Value = "'  +4.78%  "
$ws.Range('D34').Value = "'3.620"
$ws.Range('E34').Value = "'  +1.12%  "
$ws.Range('D35').Value = "'0.02227"
$ws.Range('E35').Value = "'  +1.60%  "
$ws.Range('E36').Value = "'  +1.65%  "
$ws.Range('D37').Value = "'11.20"
$ws.Range('E37').Value = "'  +1.28%  "
$ws.Range('D38').Value = "'4.839"
$ws.Range('E38').Value = "'  +1.91%  "
$ws.Range('B39').Value = 'TheSandbox'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D39').Value = "'0.6129"
$ws.Range('E39').Value = "'  +3.92%  "
$ws.Range('B40').Value = 'Algorand'
$ws.Range('C40').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D40').Value = "'0.2020"
$ws.Range('E40').Value = "'  +1.34%  "
$ws.Range('D41').Value = "'1.432"
$ws.Range('E41').Value = "'  +2.01%  "
$ws.Range('D42').Value = "'8.065"
$ws.Range('E42').Value = "'  +8.27%  "
$ws.Range('D43').Value = "'1.138"
$ws.Range('E43').Value = "'  +2.64%  "
$ws.Range('D44').Value = "'13.06"
$ws.Range('E44').Value = "'  +3.18%  "
$ws.Range('D45').Value = "'3.627"
$ws.Range('E45').Value = "'  +1.95%  "
$ws.Range('D46').Value = "'0.5744"
$ws.Range('E46').Value = "'  +2.23%  "
$ws.Range('D47').Value = "'120.81"
$ws.Range('E47').Value = "'  +2.55%  "
$ws.Range('D48').Value = "'1.888"
$ws.Range('D49').Value = "'1.112"
$ws.Range('E49').Value = "'  +2.44%  "
$ws.Range('D50').Value = "'0.06711"
$ws.Range('E50').Value = "'  +1.33%  "
$ws.Range('D51').Value = "'70.49"
$ws.Range('E51').Value = "'  +1.15%  "
